$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Match the header formatting already used by the other header cells (e.g. L1):
# bold font, thin border on all sides, centered horizontally, top-aligned vertically.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20140882
    $ws.Cells.Item($r, 15).Value = 0
}
